$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$rPrXml = '<w:rPr><w:bCs/><w:color w:val="FF0000"/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr>'

$normalParaXml = '<w:p ' + $wNs + '><w:pPr>' + $rPrXml + '</w:pPr><w:r>' + $rPrXml + '<w:t>Changes made in Experimental Branch</w:t></w:r></w:p>'

$lastParaXml = '<w:p ' + $wNs + '><w:pPr>' + $rPrXml + '</w:pPr><w:r>' + $rPrXml + '<w:lastRenderedPageBreak/><w:t>Changes made in Experimental Branch</w:t></w:r></w:p>'

$totalNew = 16
for ($i = 1; $i -le $totalNew; $i++) {
    $rng = $d.Content
    $rng.Collapse(0)
    if ($i -eq $totalNew) {
        $rng.InsertXML($lastParaXml) | Out-Null
    } else {
        $rng.InsertXML($normalParaXml) | Out-Null
    }
}
